$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): insert "Data" at C1, "Confirmado" at D1, move "Quantidade" to E1 ---
$ws.Range("C1").Value = "Data"
$ws.Range("D1").Value = "Confirmado"
$ws.Range("E1").Value = "Quantidade"
# Header style (bold, border, centered) carried over from the old C1/D1 cells
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 2 ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "Não"
$ws.Range("D2").Value = "Não"
$ws.Range("E2").Value = 0

# --- Row 3 (kept mostly blank, just shifted; C3/D3 were already blank placeholders
#     and stay untouched, E3 is a new blank placeholder cell) ---
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("E3").Formula = '=""'

# --- Row 4 ---
$ws.Range("A4").Value = 992575078
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "Não"
$ws.Range("D4").Value = "Não"
$ws.Range("E4").Value = 0

# --- Row 5 (new data, replaces old row 4's tail / shifts old row 5 data down) ---
$ws.Range("A5").Value = 990033942
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = "2023-10-19T16:45"
$ws.Range("D5").Value = "Não"
$ws.Range("E5").Value = 3

# --- Row 6 (new row, carries the old row 5 phone number with new values) ---
$ws.Range("A6").Value = 966652864
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "Não"
$ws.Range("D6").Value = "Não"
$ws.Range("E6").Value = 0
